# Add a new "Pseudonym" protocol field to the Person sheet, right after
# "Name" (i.e. before the existing "Gender" row), shifting every row
# below it down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Person")

# --- Insert the new row and fill it in -------------------------------------
# Row 3 currently holds "Gender"; inserting here pushes Gender (and every
# row after it) down by one and the new blank row inherits that row's
# formatting (matches the diff: same styles as the old row 3).
$ws.Rows.Item(3).Insert() | Out-Null

$ws.Range("A3").Value = "Person"
$ws.Range("B3").Value = "Pseudonym"
$ws.Range("C3").Value = "Enter one or more pseudonyms seperated by a comma e.g. George Sand, Aurore"
$ws.Rows.Item(3).RowHeight = 16

# --- Re-anchor the two threaded comments that sat on the now-shifted rows --
# They were on A14 ("Image" field) and B15 ("Permission" field); those
# fields are now on row 15 and row 16 respectively. The COM model has no
# "move" operation for a comment, so recreate it at the new cell using the
# the same text and delete the old one.
$ct1 = $ws.Range("A14").CommentThreaded()
$ct1Text = $ct1.Text()
$ct1.Delete() | Out-Null
$ws.Range("A15").AddCommentThreaded($ct1Text) | Out-Null

$ct2 = $ws.Range("B15").CommentThreaded()
$ct2Text = $ct2.Text()
$ct2.Delete() | Out-Null
$ws.Range("B16").AddCommentThreaded($ct2Text) | Out-Null

# --- Make "Person" the active/selected tab (it was "Picturestory") ---------
$ws.Range("C3").Select() | Out-Null
